# ------------------------------------------------------------------
# Adds a "Player Info" sheet and an "ODI Batting Extra" sheet around the
# existing "ODI Batting" sheet, and reworks the MATCH_CARD_LINK column
# (full scorecard URL) into a bare MATCH_CODE column on "ODI Batting".
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# The workbook currently holds a single sheet: "ODI Batting" with all of
# the per-match batting rows.
$wsInfo = $wb.Worksheets.Item(1)

# Duplicate it immediately after itself -- the copy keeps every row/value
# (and shared style) so we only have to touch the MATCH_CARD_LINK column
# instead of retyping 94 rows of data. Rename the original out of the way
# first so the two sheets never collide on the "ODI Batting" name.
$wsInfo.Name = "Player Info Staging"
$wsInfo.Copy([Type]::Missing, $wsInfo)
$wsBatting = $wb.Worksheets.Item(2)
$wsBatting.Name = "ODI Batting"

# Add the (initially empty) "ODI Batting Extra" sheet right after it.
$wsBatting.Copy([Type]::Missing, $wsBatting)
$wsExtra = $wb.Worksheets.Item(3)
$wsExtra.Name = "ODI Batting Extra"
$wsExtra.Cells.Clear()

# --- 1. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, and shrink
#        every row's full scorecard URL down to the bare numeric code.
$wsBatting.Range("D1").Value = "MATCH_CODE"

$lastRow = $wsBatting.UsedRange.Rows.Count
$wsBatting.Range("D2:D$lastRow").NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $txt = $cell.Text
    if ($txt -match 'MatchCode=(\d+)') {
        $cell.Value = $matches[1]
    }
}

# A handful of rows never had an INNING_NUMBER recorded; they carried an
# explicit-but-empty cell before, now they should be genuinely blank.
$blankInningRows = @(3,4,8,11,17,20,31,37,57,66,67,79,85,90,94)
foreach ($r in $blankInningRows) {
    $wsBatting.Cells.Item($r, 2).ClearContents()
}

# --- 2. "ODI Batting Extra": per-match extra batting detail.
$extraData = @(
    @('MATCH_CODE','BATTING_POSITION','NUM_4','NUM_6','PERCENT_RUNS_OF_TOTAL','MAN_OF_MATCH'),
    @('4088',$null,$null,$null,$null,'NO'),
    @('4089',7,'0','0','1.19%','NO'),
    @('4096',4,'0','0',$null,'NO'),
    @('4098',$null,$null,$null,$null,'NO'),
    @('4099',4,'3','0','11.87%','NO'),
    @('4175',$null,$null,$null,$null,'NO'),
    @('4196',$null,$null,$null,$null,'NO'),
    @('4197',4,'2','1','18.90%','NO'),
    @('4199',5,'0','0','0.57%','NO'),
    @('4201',6,$null,$null,$null,'NO'),
    @('4203',$null,$null,$null,$null,'NO'),
    @('4205',4,'1','1','16.59%','NO'),
    @('4234',6,'0','0','4.72%','NO'),
    @('4235',$null,$null,$null,$null,'NO'),
    @('4236',6,$null,$null,$null,'NO'),
    @('4245',5,'5','1','15.51%','NO'),
    @('4248',5,'0','0',$null,'NO'),
    @('4345',7,'1','0','2.55%','NO'),
    @('4350',7,$null,$null,$null,'NO'),
    @('4353',$null,$null,$null,$null,'NO')
)

# Header row + text columns (A, C, D, E, F) all need NumberFormat "@" up
# front so numeric-looking strings ("0", "1.19%", match codes, ...) land
# as text instead of getting auto-coerced to numbers. Column B
# (BATTING_POSITION) is a genuine number, so it's left as General and
# assigned bare int literals from $extraData.
$extraRowCount = $extraData.Count
$wsExtra.Range("A1:A$extraRowCount").NumberFormat = "@"
$wsExtra.Range("C1:F$extraRowCount").NumberFormat = "@"

for ($i = 0; $i -lt $extraData.Count; $i++) {
    $r = $i + 1
    $rowVals = $extraData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $val = $rowVals[$c - 1]
        $cell = $wsExtra.Cells.Item($r, $c)
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# --- 3. Turn the original sheet into "Player Info" (one row about the
#        player this workbook is for).
$wsInfo.Name = "Player Info"
$wsInfo.Cells.Clear()


$wsInfo.Range("A1:D2").NumberFormat = "@"
$playerInfoData = @(
    @('ID','NAME','BATTING_HAND','BOWL_STYLE'),
    @('3210','Krishnakumar Dinesh Karthik','Right Handed','Does Not Bowl | Unknown')
)
for ($i = 0; $i -lt $playerInfoData.Count; $i++) {
    $r = $i + 1
    $rowVals = $playerInfoData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $wsInfo.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
